$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-10-29 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-30 Wednesday", 2) | Out-Null

# Update each math-problem cell in the table (20 rows x 5 cols)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "75-68="
$t.Cell(1,2).Range.Text = "43+48="
$t.Cell(1,3).Range.Text = "47+27="
$t.Cell(1,4).Range.Text = "24+19="
$t.Cell(1,5).Range.Text = "31-7="
$t.Cell(2,1).Range.Text = "40-13="
$t.Cell(2,2).Range.Text = "7+24="
$t.Cell(2,3).Range.Text = "6+38="
$t.Cell(2,4).Range.Text = "14+29="
$t.Cell(2,5).Range.Text = "3+29="
$t.Cell(3,1).Range.Text = "55+38="
$t.Cell(3,2).Range.Text = "97-59="
$t.Cell(3,3).Range.Text = "66+6="
$t.Cell(3,4).Range.Text = "39+59="
$t.Cell(3,5).Range.Text = "85-18="
$t.Cell(4,1).Range.Text = "27+7="
$t.Cell(4,2).Range.Text = "3+28="
$t.Cell(4,3).Range.Text = "56+28="
$t.Cell(4,4).Range.Text = "72-53="
$t.Cell(4,5).Range.Text = "85-19="
$t.Cell(5,1).Range.Text = "71-42="
$t.Cell(5,2).Range.Text = "64-39="
$t.Cell(5,3).Range.Text = "24+38="
$t.Cell(5,4).Range.Text = "7+84="
$t.Cell(5,5).Range.Text = "36+46="
$t.Cell(6,1).Range.Text = "28+56="
$t.Cell(6,2).Range.Text = "81-44="
$t.Cell(6,3).Range.Text = "14+68="
$t.Cell(6,4).Range.Text = "34+27="
$t.Cell(6,5).Range.Text = "7+39="
$t.Cell(7,1).Range.Text = "62-57="
$t.Cell(7,2).Range.Text = "83-67="
$t.Cell(7,3).Range.Text = "36+38="
$t.Cell(7,4).Range.Text = "82-48="
$t.Cell(7,5).Range.Text = "78+4="
$t.Cell(8,1).Range.Text = "83-65="
$t.Cell(8,2).Range.Text = "41-13="
$t.Cell(8,3).Range.Text = "80-31="
$t.Cell(8,4).Range.Text = "27+39="
$t.Cell(8,5).Range.Text = "82-75="
$t.Cell(9,1).Range.Text = "54+8="
$t.Cell(9,2).Range.Text = "16+77="
$t.Cell(9,3).Range.Text = "80-45="
$t.Cell(9,4).Range.Text = "29+14="
$t.Cell(9,5).Range.Text = "44-37="
$t.Cell(10,1).Range.Text = "36+45="
$t.Cell(10,2).Range.Text = "84-27="
$t.Cell(10,3).Range.Text = "85-69="
$t.Cell(10,4).Range.Text = "9+18="
$t.Cell(10,5).Range.Text = "35+28="
$t.Cell(11,1).Range.Text = "30-12="
$t.Cell(11,2).Range.Text = "16+45="
$t.Cell(11,3).Range.Text = "95-59="
$t.Cell(11,4).Range.Text = "81-8="
$t.Cell(11,5).Range.Text = "45+26="
$t.Cell(12,1).Range.Text = "61-22="
$t.Cell(12,2).Range.Text = "29+22="
$t.Cell(12,3).Range.Text = "40-37="
$t.Cell(12,4).Range.Text = "76-9="
$t.Cell(12,5).Range.Text = "24+27="
$t.Cell(13,1).Range.Text = "52-14="
$t.Cell(13,2).Range.Text = "86-38="
$t.Cell(13,3).Range.Text = "61-37="
$t.Cell(13,4).Range.Text = "93-76="
$t.Cell(13,5).Range.Text = "49+14="
$t.Cell(14,1).Range.Text = "91-44="
$t.Cell(14,2).Range.Text = "95-17="
$t.Cell(14,3).Range.Text = "40-32="
$t.Cell(14,4).Range.Text = "42-13="
$t.Cell(14,5).Range.Text = "28+23="
$t.Cell(15,1).Range.Text = "71-49="
$t.Cell(15,2).Range.Text = "6+85="
$t.Cell(15,3).Range.Text = "62-24="
$t.Cell(15,4).Range.Text = "90-4="
$t.Cell(15,5).Range.Text = "45-28="
$t.Cell(16,1).Range.Text = "74-7="
$t.Cell(16,2).Range.Text = "27+6="
$t.Cell(16,3).Range.Text = "8+19="
$t.Cell(16,4).Range.Text = "5+69="
$t.Cell(16,5).Range.Text = "24+18="
$t.Cell(17,1).Range.Text = "80-33="
$t.Cell(17,2).Range.Text = "9+64="
$t.Cell(17,3).Range.Text = "17+38="
$t.Cell(17,4).Range.Text = "67+6="
$t.Cell(17,5).Range.Text = "80-38="
$t.Cell(18,1).Range.Text = "92-86="
$t.Cell(18,2).Range.Text = "93-56="
$t.Cell(18,3).Range.Text = "38+6="
$t.Cell(18,4).Range.Text = "40-12="
$t.Cell(18,5).Range.Text = "26+5="
$t.Cell(19,1).Range.Text = "46+37="
$t.Cell(19,2).Range.Text = "80-79="
$t.Cell(19,3).Range.Text = "57+35="
$t.Cell(19,4).Range.Text = "70-37="
$t.Cell(19,5).Range.Text = "18+4="
$t.Cell(20,1).Range.Text = "49+3="
$t.Cell(20,2).Range.Text = "26+25="
$t.Cell(20,3).Range.Text = "91-13="
$t.Cell(20,4).Range.Text = "45+46="
$t.Cell(20,5).Range.Text = "4+88="
